# Updated cryptos list on Sat May 27 04:41:37 UTC 2023 with GitHub Actions
#
# Refresh Price (D) and Volume(1h) (E) figures for every coin row, and
# swap three pairs of adjacent rows whose rank order changed so that
# Coin (B) / Link (C) also need to move:
#   rows 20/21: WrappedBTC <-> Avalanche
#   rows 34/35: Filecoin <-> RenderToken
#   rows 37/38: VeChain <-> Hedera

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: temporarily force column D to Text so numeric-looking price
# strings ("1.007", "309.10", ...) are not auto-coerced to numbers by COM,
# matching the inline-string <is><t> cells in the source file.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.914.21"
$ws.Range("E2").Value = "  +1.07%  "
$ws.Range("D3").Value = "1.845.26"
$ws.Range("E3").Value = "  +1.27%  "
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "309.10"
$ws.Range("E5").Value = "  +0.98%  "
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("D7").Value = "0.4717"
$ws.Range("E7").Value = "  +0.65%  "
$ws.Range("D8").Value = "0.3682"
$ws.Range("E8").Value = "  +2.28%  "
$ws.Range("D9").Value = "0.07218"
$ws.Range("E9").Value = "  +1.15%  "
$ws.Range("D10").Value = "0.9250"
$ws.Range("E10").Value = "  +2.43%  "
$ws.Range("D11").Value = "19.64"
$ws.Range("E11").Value = "  +0.99%  "
$ws.Range("D12").Value = "0.07621"
$ws.Range("E12").Value = "  -2.48%  "
$ws.Range("D13").Value = "1.845.82"
$ws.Range("E13").Value = "  +2.02%  "
$ws.Range("D14").Value = "5.311"
$ws.Range("E14").Value = "  +0.99%  "
$ws.Range("D15").Value = "6.400"
$ws.Range("E15").Value = "  +1.10%  "
$ws.Range("D16").Value = "88.42"
$ws.Range("E16").Value = "  +1.05%  "
$ws.Range("D17").Value = "1.008"
$ws.Range("E17").Value = "  -0.13%  "
$ws.Range("D18").Value = "0.000008665"
$ws.Range("E18").Value = "  +1.01%  "
$ws.Range("E19").Value = "  -0.21%  "
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "14.57"
$ws.Range("E20").Value = "  +2.68%  "
$ws.Range("B21").Value = "WrappedBTC"
$ws.Range("C21").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D21").Value = "26.937.45"
$ws.Range("E21").Value = "  +1.02%  "
$ws.Range("E22").Value = "  +0.33%  "
$ws.Range("E23").Value = "  +0.82%  "
$ws.Range("D24").Value = "1.920"
$ws.Range("E24").Value = "  -0.65%  "
$ws.Range("D25").Value = "152.10"
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("D26").Value = "18.14"
$ws.Range("E26").Value = "  +1.24%  "
$ws.Range("D27").Value = "2.006"
$ws.Range("E27").Value = "  +1.25%  "
$ws.Range("E28").Value = "  +0.42%  "
$ws.Range("D29").Value = "4.956"
$ws.Range("E29").Value = "  +2.90%  "
$ws.Range("D30").Value = "0.08837"
$ws.Range("E30").Value = "  +0.36%  "
$ws.Range("D31").Value = "3.304"
$ws.Range("E31").Value = "  +5.03%  "
$ws.Range("D32").Value = "0.7451"
$ws.Range("E32").Value = "  +1.84%  "
$ws.Range("D33").Value = "1.167"
$ws.Range("E33").Value = "  +3.83%  "
$ws.Range("B34").Value = "RenderToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D34").Value = "2.770"
$ws.Range("E34").Value = "  -0.18%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").Value = "4.490"
$ws.Range("E35").Value = "  +0.94%  "
$ws.Range("D36").Value = "1.090"
$ws.Range("E36").Value = "  +1.22%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "0.05261"
$ws.Range("E37").Value = "  +2.84%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.01948"
$ws.Range("E38").Value = "  +1.03%  "
$ws.Range("D39").Value = "2.961"
$ws.Range("E39").Value = "  +1.47%  "
$ws.Range("D40").Value = "0.5217"
$ws.Range("E40").Value = "  +2.90%  "
$ws.Range("D41").Value = "6.914"
$ws.Range("E41").Value = "  +1.37%  "
$ws.Range("D42").Value = "0.1513"
$ws.Range("E42").Value = "  +0.95%  "
$ws.Range("D43").Value = "8.217"
$ws.Range("E43").Value = "  +2.64%  "
$ws.Range("D44").Value = "10.56"
$ws.Range("E44").Value = "  +5.93%  "
$ws.Range("D45").Value = "0.4696"
$ws.Range("E45").Value = "  +0.29%  "
$ws.Range("D46").Value = "1.007"
$ws.Range("E46").Value = "  -0.20%  "
$ws.Range("D47").Value = "102.08"
$ws.Range("E47").Value = "  +3.12%  "
$ws.Range("D48").Value = "1.603"
$ws.Range("E48").Value = "  +2.62%  "
$ws.Range("D49").Value = "65.52"
$ws.Range("E49").Value = "  +2.92%  "
$ws.Range("D50").Value = "0.06034"
$ws.Range("E50").Value = "  +0.40%  "
$ws.Range("D51").Value = "0.8854"
$ws.Range("E51").Value = "  +4.00%  "

# Step 2: drop the temporary Text format so the cells end up with no
# explicit style override again (same as every other B:E data cell).
$ws.Range("D2:D51").ClearFormats()
